# The "day" sheet's bsecode column (D) was stored as text for rows 144-152;
# convert those nine cells to real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bsecodes = @{
    144 = 534091
    145 = 532978
    146 = 500043
    147 = 532454
    148 = 524494
    149 = 532555
    150 = 530005
    151 = 540065
    152 = 532822
}

foreach ($row in $bsecodes.Keys) {
    $ws.Range("D$row").Value = $bsecodes[$row]
}

# Append five new stock rows (153-157) that were broken out of stock.yaml.
$newRows = @(
    @(153, 1, "SHREECEM",   "Shree Cements Limited",                         "500387", -0.08, 27643.6, 65384,   "day", "15/07/2024 11:35:33"),
    @(154, 2, "PVRINOX",    "PVR Inox Ltd",                                  "532689", -1.15, 1459.95, 492900,  "day", "15/07/2024 11:35:33"),
    @(155, 3, "CANFINHOME", "Can Fin Homes Limited",                         "511196", 1.17,  878.45,  737699,  "day", "15/07/2024 11:35:33"),
    @(156, 4, "LICHSGFIN",  "Lic Housing Finance Limited",                   "500253", 3.03,  803.3,   3307875, "day", "15/07/2024 11:35:33"),
    @(157, 5, "CROMPTON",   "Crompton Greaves Consumer Electricals Limited", "539876", -0.05, 431.55,  928231,  "day", "15/07/2024 11:35:33")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    # bsecode keeps its text formatting for these newly-added rows (unlike
    # the D144-D152 fix above), so force the column to Text before writing
    # the numeric-looking code, otherwise Excel auto-converts it to a number.
    $ws.Cells.Item($rowNum, 4).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
}
